$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = "27.686.38"
$ws.Cells.Item(2, 5).Value = "  -2.00%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "1.758.86"
$ws.Cells.Item(3, 5).Value = "  -2.05%  "

# Row 4
$cell = $ws.Cells.Item(4, 4)
$cell.Formula = "=""1.008"""
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Cells.Item(4, 5).Value = "  +0.48%  "

# Row 5
$cell = $ws.Cells.Item(5, 4)
$cell.Formula = "=""326.02"""
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Cells.Item(5, 5).Value = "  -0.39%  "

# Row 6
$cell = $ws.Cells.Item(6, 4)
$cell.Formula = "=""1.003"""
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Cells.Item(6, 5).Value = "  +0.24%  "

# Row 7
$cell = $ws.Cells.Item(7, 4)
$cell.Formula = "=""0.4416"""
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Cells.Item(7, 5).Value = "  -2.39%  "

# Row 8
$cell = $ws.Cells.Item(8, 4)
$cell.Formula = "=""0.3726"""
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Cells.Item(8, 5).Value = "  -0.51%  "

# Row 9
$cell = $ws.Cells.Item(9, 4)
$cell.Formula = "=""45.99"""
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Cells.Item(9, 5).Value = "  +2.92%  "

# Row 10
$cell = $ws.Cells.Item(10, 4)
$cell.Formula = "=""0.07746"""
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Cells.Item(10, 5).Value = "  +2.61%  "

# Row 11
$cell = $ws.Cells.Item(11, 4)
$cell.Formula = "=""1.126"""
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Cells.Item(11, 5).Value = "  -2.00%  "

# Row 12
$cell = $ws.Cells.Item(12, 4)
$cell.Formula = "=""1.005"""
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Cells.Item(12, 5).Value = "  +0.36%  "

# Row 13
$cell = $ws.Cells.Item(13, 4)
$cell.Formula = "=""21.75"""
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Cells.Item(13, 5).Value = "  -3.57%  "

# Row 14
$cell = $ws.Cells.Item(14, 4)
$cell.Formula = "=""6.194"""
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Cells.Item(14, 5).Value = "  -1.75%  "

# Row 15
$cell = $ws.Cells.Item(15, 4)
$cell.Formula = "=""7.360"""
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Cells.Item(15, 5).Value = "  -2.51%  "

# Row 16
$ws.Cells.Item(16, 4).Value = "1.760.29"
$ws.Cells.Item(16, 5).Value = "  -0.03%  "

# Row 17
$cell = $ws.Cells.Item(17, 4)
$cell.Formula = "=""91.50"""
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Cells.Item(17, 5).Value = "  +12.72%  "

# Row 18
$cell = $ws.Cells.Item(18, 4)
$cell.Formula = "=""0.00001080"""
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Cells.Item(18, 5).Value = "  -1.11%  "

# Row 19
$ws.Cells.Item(19, 5).Value = "  -7.47%  "

# Row 20
$ws.Cells.Item(20, 5).Value = "  +0.17%  "

# Row 21
$cell = $ws.Cells.Item(21, 4)
$cell.Formula = "=""17.36"""
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Cells.Item(21, 5).Value = "  -1.24%  "

# Row 22
$cell = $ws.Cells.Item(22, 4)
$cell.Formula = "=""6.194"""
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Cells.Item(22, 5).Value = "  -2.59%  "

# Row 23
$cell = $ws.Cells.Item(23, 4)
$cell.Formula = "=""0.5324"""
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Cells.Item(23, 5).Value = "  -1.38%  "

# Row 24
$ws.Cells.Item(24, 4).Value = "27.740.18"
$ws.Cells.Item(24, 5).Value = "  -1.74%  "

# Row 25
$cell = $ws.Cells.Item(25, 4)
$cell.Formula = "=""11.65"""
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Cells.Item(25, 5).Value = "  -1.21%  "

# Row 26
$cell = $ws.Cells.Item(26, 4)
$cell.Formula = "=""2.334"""
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Cells.Item(26, 5).Value = "  -3.73%  "

# Row 27
$cell = $ws.Cells.Item(27, 4)
$cell.Formula = "=""20.81"""
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Cells.Item(27, 5).Value = "  +1.19%  "

# Row 28
$cell = $ws.Cells.Item(28, 4)
$cell.Formula = "=""153.76"""
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Cells.Item(28, 5).Value = "  +1.50%  "

# Row 29
$cell = $ws.Cells.Item(29, 4)
$cell.Formula = "=""2.366"""
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Cells.Item(29, 5).Value = "  +0.31%  "

# Row 30
$ws.Cells.Item(30, 4).Value = "1.959.12"
$ws.Cells.Item(30, 5).Value = "  -1.88%  "

# Row 31
$cell = $ws.Cells.Item(31, 4)
$cell.Formula = "=""129.08"""
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Cells.Item(31, 5).Value = "  -3.18%  "

# Row 32
$cell = $ws.Cells.Item(32, 4)
$cell.Formula = "=""1.210"""
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Cells.Item(32, 5).Value = "  -2.12%  "

# Row 33
$cell = $ws.Cells.Item(33, 4)
$cell.Formula = "=""5.767"""
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Cells.Item(33, 5).Value = "  -1.05%  "

# Row 34
$cell = $ws.Cells.Item(34, 4)
$cell.Formula = "=""0.09274"""
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Cells.Item(34, 5).Value = "  -1.78%  "

# Row 35
$cell = $ws.Cells.Item(35, 4)
$cell.Formula = "=""3.682"""
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Cells.Item(35, 5).Value = "  -8.53%  "

# Row 36
$cell = $ws.Cells.Item(36, 4)
$cell.Formula = "=""12.76"""
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Cells.Item(36, 5).Value = "  +5.00%  "

# Row 37
$ws.Cells.Item(37, 2).Value = "Algorand"
$ws.Cells.Item(37, 3).Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$cell = $ws.Cells.Item(37, 4)
$cell.Formula = "=""0.2191"""
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Cells.Item(37, 5).Value = "  -6.35%  "

# Row 38
$ws.Cells.Item(38, 2).Value = "VeChain"
$ws.Cells.Item(38, 3).Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$cell = $ws.Cells.Item(38, 4)
$cell.Formula = "=""0.02335"""
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Cells.Item(38, 5).Value = "  -0.01%  "

# Row 39
$cell = $ws.Cells.Item(39, 4)
$cell.Formula = "=""0.6513"""
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Cells.Item(39, 5).Value = "  -0.96%  "

# Row 40
$cell = $ws.Cells.Item(40, 4)
$cell.Formula = "=""5.101"""
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Cells.Item(40, 5).Value = "  -1.58%  "

# Row 41
$cell = $ws.Cells.Item(41, 4)
$cell.Formula = "=""0.06129"""
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Cells.Item(41, 5).Value = "  -3.50%  "

# Row 42
$cell = $ws.Cells.Item(42, 4)
$cell.Formula = "=""1.197"""
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Cells.Item(42, 5).Value = "  -0.89%  "

# Row 43
$cell = $ws.Cells.Item(43, 4)
$cell.Formula = "=""8.008"""
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Cells.Item(43, 5).Value = "  -4.32%  "

# Row 44
$cell = $ws.Cells.Item(44, 4)
$cell.Formula = "=""1.415"""
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Cells.Item(44, 5).Value = "  -3.88%  "

# Row 45
$cell = $ws.Cells.Item(45, 4)
$cell.Formula = "=""1.003"""
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Cells.Item(45, 5).Value = "  +0.30%  "

# Row 46
$cell = $ws.Cells.Item(46, 4)
$cell.Formula = "=""13.92"""
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Cells.Item(46, 5).Value = "  -1.25%  "

# Row 47
$cell = $ws.Cells.Item(47, 4)
$cell.Formula = "=""0.6015"""
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Cells.Item(47, 5).Value = "  -1.44%  "

# Row 48
$cell = $ws.Cells.Item(48, 4)
$cell.Formula = "=""3.756"""
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Cells.Item(48, 5).Value = "  -1.05%  "

# Row 49
$cell = $ws.Cells.Item(49, 4)
$cell.Formula = "=""126.15"""
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Cells.Item(49, 5).Value = "  -2.82%  "

# Row 50
$cell = $ws.Cells.Item(50, 4)
$cell.Formula = "=""2.000"""
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Cells.Item(50, 5).Value = "  -1.50%  "

# Row 51
$cell = $ws.Cells.Item(51, 4)
$cell.Formula = "=""1.146"""
$cell.Copy()
$cell.PasteSpecial(-4163)
$ws.Cells.Item(51, 5).Value = "  -1.47%  "

$excel.CutCopyMode = 0
